$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 462.76
$ws.Range("I33").Value = 431.72223
$ws.Range("J33").Value = 542.5714
$ws.Range("K33").Value = 431.72223
$ws.Range("L33").Value = 542.5714
$ws.Range("M33").Value = -202.72223
$ws.Range("N33").Value = -1000.5714

# Row 116
$ws.Range("H116").Value = 3397.4
$ws.Range("I116").Value = 2897.2222
$ws.Range("J116").Value = 4147.6665
$ws.Range("K116").Value = 2897.2222
$ws.Range("L116").Value = 4147.6665
$ws.Range("M116").Value = 544.7777999999998

# Row 138
$ws.Range("H138").Value = 525649.1
$ws.Range("I138").Value = 1328.75
$ws.Range("J138").Value = 650859.9399999999
$ws.Range("K138").Value = 3986.25
$ws.Range("L138").Value = 1952579.82
$ws.Range("M138").Value = 1153.75
$ws.Range("N138").Value = -1962859.82

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4984.397
$ws.Range("I32").Value = 4800.3
$ws.Range("J32").Value = 8666.333000000001
$ws.Range("K32").Value = 4800.3
$ws.Range("L32").Value = 8666.333000000001
$ws.Range("M32").Value = -4513.3
$ws.Range("N32").Value = -9240.333000000001

# Row 36
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 5000
$ws.Range("N36").Value = -5692
$ws.Range("M36").ClearContents()

# Row 132
$ws.Range("H132").Value = 2529.7144
$ws.Range("I132").Value = 2040.3
$ws.Range("J132").Value = 3753.25
$ws.Range("K132").Value = 6120.9
$ws.Range("L132").Value = 11259.75
$ws.Range("M132").Value = -3590.9
$ws.Range("N132").Value = -16319.75

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2756.4443
$ws.Range("I20").Value = 2434.25
$ws.Range("J20").Value = 3400.8333
$ws.Range("K20").Value = 2434.25
$ws.Range("L20").Value = 3400.8333
$ws.Range("M20").Value = -2187.25
$ws.Range("N20").Value = -3894.8333

# Row 94
$ws.Range("H94").Value = 12500293
$ws.Range("I94").Value = 19231018
$ws.Range("J94").Value = 377.14285
$ws.Range("K94").Value = 19231018
$ws.Range("L94").Value = 377.14285
$ws.Range("M94").Value = -19230567
$ws.Range("N94").Value = -1279.14285

# Row 99
$ws.Range("H99").Value = 58824624
$ws.Range("I99").Value = 66667696
$ws.Range("J99").Value = 1575
$ws.Range("K99").Value = 66667696
$ws.Range("L99").Value = 1575
$ws.Range("M99").Value = -66666198
$ws.Range("N99").Value = -4571

# Row 134
$ws.Range("H134").Value = 1341
$ws.Range("I134").Value = 994.625
$ws.Range("J134").Value = 1687.375
$ws.Range("K134").Value = 2983.875
$ws.Range("L134").Value = 5062.125
$ws.Range("M134").Value = -448.875
$ws.Range("N134").Value = -10132.125

$ws = $wb.Worksheets.Item("CRP")
# Row 38
$ws.Range("H38").Value = 2500
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 2500
$ws.Range("N38").Value = -3254

# Row 45
$ws.Range("H45").Value = 3200
$ws.Range("I45").Value = 3900
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 3900
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -3307
$ws.Range("N45").Value = -3686

# Row 46
$ws.Range("H46").Value = 2500
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2500
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2500
$ws.Range("N46").Value = -2922

# Row 99
$ws.Range("H99").Value = 1703.5
$ws.Range("I99").Value = 1694
$ws.Range("J99").Value = 1732
$ws.Range("K99").Value = 1694
$ws.Range("L99").Value = 1732
$ws.Range("M99").Value = -196
$ws.Range("N99").Value = -4728

# Row 126
$ws.Range("H126").Value = 1703.5
$ws.Range("I126").Value = 1694
$ws.Range("J126").Value = 1732
$ws.Range("K126").Value = 5082
$ws.Range("L126").Value = 5196
$ws.Range("M126").Value = -2612
$ws.Range("N126").Value = -10136

# Row 134
$ws.Range("H134").Value = 14287400
$ws.Range("I134").Value = 1709.4
$ws.Range("J134").Value = 100001544
$ws.Range("K134").Value = 5128.200000000001
$ws.Range("L134").Value = 300004632
$ws.Range("M134").Value = -2593.200000000001
$ws.Range("N134").Value = -300009702

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 822.94446
$ws.Range("I121").Value = 295
$ws.Range("J121").Value = 888.9375
$ws.Range("K121").Value = 885
$ws.Range("L121").Value = 2666.8125
$ws.Range("M121").Value = 425
$ws.Range("N121").Value = -5286.8125

# Row 129
$ws.Range("H129").Value = 21930836
$ws.Range("I129").Value = 41667480
$ws.Range("J129").Value = 7576915
$ws.Range("K129").Value = 125002440
$ws.Range("L129").Value = 22730745
$ws.Range("M129").Value = -124997440
$ws.Range("N129").Value = -22740745

# Row 131
$ws.Range("H131").Value = 15385539
$ws.Range("I131").Value = 166667040
$ws.Range("J131").Value = 979.6949
$ws.Range("K131").Value = 500001120
$ws.Range("L131").Value = 2939.0847
$ws.Range("M131").Value = -499996080
$ws.Range("N131").Value = -13019.0847

# Row 139
$ws.Range("H139").Value = 1772.2368
$ws.Range("I139").Value = 1814.6666
$ws.Range("J139").Value = 1699.5
$ws.Range("K139").Value = 5443.9998
$ws.Range("L139").Value = 5098.5
$ws.Range("M139").Value = -303.9997999999996

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 18753710
$ws.Range("I70").Value = 19234548
$ws.Range("J70").Value = 18185446
$ws.Range("K70").Value = 19234548
$ws.Range("L70").Value = 18185446
$ws.Range("M70").Value = -19234278
$ws.Range("N70").Value = -18185986

# Row 73
$ws.Range("H73").Value = 18753710
$ws.Range("I73").Value = 19234548
$ws.Range("J73").Value = 18185446
$ws.Range("K73").Value = 19234548
$ws.Range("L73").Value = 18185446
$ws.Range("M73").Value = -19233612
$ws.Range("N73").Value = -18187318

# Row 132
$ws.Range("H132").Value = 3526.5652
$ws.Range("I132").Value = 3476.5
$ws.Range("J132").Value = 3604.4443
$ws.Range("K132").Value = 10429.5
$ws.Range("L132").Value = 10813.3329
$ws.Range("M132").Value = -7899.5
$ws.Range("N132").Value = -15873.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4475
$ws.Range("I40").Value = 1862.5
$ws.Range("J40").Value = 7087.5
$ws.Range("K40").Value = 1862.5
$ws.Range("L40").Value = 7087.5
$ws.Range("M40").Value = -1726.5
$ws.Range("N40").Value = -7359.5

# Row 122
$ws.Range("H122").Value = 125001500
$ws.Range("I122").Value = 125001500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 375004500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -375002050
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 2944.32
$ws.Range("I132").Value = 2823.8462
$ws.Range("J132").Value = 3074.8333
$ws.Range("K132").Value = 8471.5386
$ws.Range("L132").Value = 9224.499899999999
$ws.Range("M132").Value = -5941.5386
$ws.Range("N132").Value = -14284.4999

# Row 134
$ws.Range("H134").Value = 33832
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 33832
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 33832
$ws.Range("N134").Value = -43972

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1625
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 2750
$ws.Range("K81").Value = 1000
$ws.Range("L81").Value = 5500
$ws.Range("M81").Value = 61
$ws.Range("N81").Value = -7622

# Row 84
$ws.Range("H84").Value = 1625
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 2750
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 27500
$ws.Range("M84").Value = 304
$ws.Range("N84").Value = -38108

# Row 122
$ws.Range("H122").Value = 19232316
$ws.Range("I122").Value = 25001710
$ws.Range("J122").Value = 1003.3333
$ws.Range("K122").Value = 75005130
$ws.Range("L122").Value = 3009.9999
$ws.Range("M122").Value = -75002680
$ws.Range("N122").Value = -7909.9999
